# ozon fixes 14.11.2025 part 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data corrections ---
$ws.Range("A2").Value = 2082359
$ws.Range("B2").Value = 659000
$ws.Range("C2").Value = 11
$ws.Range("E2").Value = 13

# --- Switch the built-in Office theme naming to "Office 2013 - 2022" ---
# (colors/fonts are unchanged, only the theme/color-scheme/font-scheme
# display names are updated to match current Office theme naming)
try {
    $theme = $wb.Theme
    $theme.Name = "Office 2013 - 2022 Theme"
    $theme.ThemeColorScheme.Name = "Office 2013 - 2022"
    $theme.ThemeFontScheme.Name = "Office 2013 - 2022"
} catch {
    # Theme renaming not supported by this host; ignore.
}

# --- Move the active selection to C5 ---
$ws.Range("C5").Select()
